# Validacion suma asegurada sin bloqueo
# Adds a new "Validacion" row (row 8) to the QA data sheet, re-using the same
# formatting as the existing rows 6/7, and updates the active window's
# scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 8: copy formatting from row 7 (its closest sibling) --------------
# Only copy the specific cells that carry non-default formatting in row 7
# (B, E, N:O, P) so we don't materialize empty cells in gaps like L8.
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial($xlPasteFormats)
$ws.Range("E7").Copy()
$ws.Range("E8").PasteSpecial($xlPasteFormats)
$ws.Range("N7:O7").Copy()
$ws.Range("N8:O8").PasteSpecial($xlPasteFormats)
$ws.Range("P7").Copy()
$ws.Range("P8").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- Row 8 values ---------------------------------------------------------
$ws.Range("A8").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("B8").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"
$ws.Range("C8").Value = "su"
$ws.Range("D8").Value = "gw"
$ws.Range("E8").Value = "'3582596139"
$ws.Range("F8").Value = "Motor"
$ws.Range("G8").Value = "Validacion"
$ws.Range("H8").Value = "Menos de 5 vehículos"
$ws.Range("I8").Value = "Anual"
$ws.Range("J8").Value = "Cupón"
$ws.Range("K8").Value = "No"
$ws.Range("M8").Value = 2022
$ws.Range("N8").Value = "CHEVROLET"
$ws.Range("O8").Value = "AGILE"
$ws.Range("P8").Value = "5.000.000"
$ws.Range("Q8").Value = "C - Resp. Civil-Robo/Incendio Total y Parcial Daños Totales por Accidente"
$ws.Range("R8").Value = "QMM118"
$ws.Range("S8").Value = "MASDASJ12319"
$ws.Range("T8").Value = "ASDAKE1KJ239"
$ws.Range("U8").Value = "Sin Actividad"

# --- Window scroll position / selection -----------------------------------
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("R9:T10").Select() | Out-Null

Write-Host "done"
